$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 7.809785945681368
$ws.Range("D2").Value = 9.595176924504266
$ws.Range("E2").Value = 10.48731393807477
$ws.Range("F2").Value = 67.70479850439743
$ws.Range("G2").Value = 3.827234268252599
$ws.Range("J2").Value = 8.702682834169476
$ws.Range("M2").Value = 48.62621894029864

$ws.Range("B3").Value = 7.739521038045326
$ws.Range("D3").Value = 9.286498472122851
$ws.Range("E3").Value = 9.954504261527433
$ws.Range("F3").Value = 67.90224288708461
$ws.Range("G3").Value = 3.840385504923581
$ws.Range("J3").Value = 8.701777627167047
$ws.Range("M3").Value = 46.90597461090255

$ws.Range("B4").Value = 7.697950958231699
$ws.Range("D4").Value = 9.098567622237313
$ws.Range("E4").Value = 9.612656639318635
$ws.Range("F4").Value = 68.07359933462178
$ws.Range("G4").Value = 3.848790165778024
$ws.Range("J4").Value = 8.703117380861245
$ws.Range("M4").Value = 45.82440484611626

$ws.Range("B5").Value = 7.681423838819826
$ws.Range("D5").Value = 9.022505779376933
$ws.Range("E5").Value = 9.469748082937125
$ws.Range("F5").Value = 68.15569381798305
$ws.Range("G5").Value = 3.852299102063601
$ws.Range("J5").Value = 8.70413748162049
$ws.Range("M5").Value = 45.37783001061641

$ws.Range("B6").Value = 7.678704966400736
$ws.Range("D6").Value = 9.00991083914751
$ws.Range("E6").Value = 9.445803743970353
$ws.Range("F6").Value = 68.17005635129567
$ws.Range("G6").Value = 3.852886859969758
$ws.Range("J6").Value = 8.704335432736256
$ws.Range("M6").Value = 45.30334169649888

$ws.Range("B7").Value = 7.697726372783553
$ws.Range("D7").Value = 9.097539553028598
$ws.Range("E7").Value = 9.610743774440101
$ws.Range("F7").Value = 68.07465727833184
$ws.Range("G7").Value = 3.84883714716115
$ws.Range("J7").Value = 8.70312922179996
$ws.Range("M7").Value = 45.81840501256553

$ws.Range("B8").Value = 7.785242885201263
$ws.Range("D8").Value = 9.488482855354874
$ws.Range("E8").Value = 10.30670252782794
$ws.Range("F8").Value = 67.7622714834104
$ws.Range("G8").Value = 3.831701018102008
$ws.Range("J8").Value = 8.701975764869605
$ws.Range("M8").Value = 48.03864523142493

$ws.Range("B9").Value = 7.968540703420469
$ws.Range("D9").Value = 10.2628575485039
$ws.Range("E9").Value = 11.55164613679425
$ws.Range("F9").Value = 67.56227633002923
$ws.Range("G9").Value = 3.800664302188297
$ws.Range("J9").Value = 8.714873064270868
$ws.Range("M9").Value = 52.17053059572703

$ws.Range("B10").Value = 8.10921880215345
$ws.Range("D10").Value = 10.83042408048767
$ws.Range("E10").Value = 12.3900921466708
$ws.Range("F10").Value = 67.68668823783406
$ws.Range("G10").Value = 3.779355562539251
$ws.Range("J10").Value = 8.733763698564717
$ws.Range("M10").Value = 55.04673003044748

$ws.Range("B11").Value = 8.17427845405293
$ws.Range("D11").Value = 11.08722134229435
$ws.Range("E11").Value = 12.75456755066412
$ws.Range("F11").Value = 67.80626063797041
$ws.Range("G11").Value = 3.769969847501043
$ws.Range("J11").Value = 8.744441251372963
$ws.Range("M11").Value = 56.31650898230983

$ws.Range("B12").Value = 8.199047286595276
$ws.Range("D12").Value = 11.18418456959567
$ws.Range("E12").Value = 12.8901284558734
$ws.Range("F12").Value = 67.8609116858056
$ws.Range("G12").Value = 3.766458598354824
$ws.Range("J12").Value = 8.748787516804226
$ws.Range("M12").Value = 56.79151453017241

$ws.Range("B13").Value = 8.193707330326975
$ws.Range("D13").Value = 11.16331543901808
$ws.Range("E13").Value = 12.86104257421104
$ws.Range("F13").Value = 67.84871991474056
$ws.Range("G13").Value = 3.76721292181091
$ws.Range("J13").Value = 8.747837946267714
$ws.Range("M13").Value = 56.68947694055954

$ws.Range("B14").Value = 8.176313671913036
$ws.Range("D14").Value = 11.09520460637958
$ws.Range("E14").Value = 12.7657695873896
$ws.Range("F14").Value = 67.81056723060021
$ws.Range("G14").Value = 3.769680122648011
$ws.Range("J14").Value = 8.744792731635254
$ws.Range("M14").Value = 56.35570629872441

$ws.Range("B15").Value = 8.16567614408944
$ws.Range("D15").Value = 11.05344604302815
$ws.Range("E15").Value = 12.70709151115613
$ws.Range("F15").Value = 67.78842721186949
$ws.Range("G15").Value = 3.771196902892673
$ws.Range("J15").Value = 8.742966994955561
$ws.Range("M15").Value = 56.15049510918878

$ws.Range("B16").Value = 8.104986632606087
$ws.Range("D16").Value = 10.81360668235749
$ws.Range("E16").Value = 12.3659298717299
$ws.Range("F16").Value = 67.68016842258555
$ws.Range("G16").Value = 3.779975020784192
$ws.Range("J16").Value = 8.733108104958955
$ws.Range("M16").Value = 54.96294501323398

$ws.Range("B17").Value = 8.068013339732488
$ws.Range("D17").Value = 10.66605626588436
$ws.Range("E17").Value = 12.15228126010263
$ws.Range("F17").Value = 67.63011356985629
$ws.Range("G17").Value = 3.785437999307208
$ws.Range("J17").Value = 8.727596033983961
$ws.Range("M17").Value = 54.22431746391823

$ws.Range("B18").Value = 8.046848728545744
$ws.Range("D18").Value = 10.5810623388357
$ws.Range("E18").Value = 12.02780267754559
$ws.Range("F18").Value = 67.60723290710355
$ws.Range("G18").Value = 3.788609209061694
$ws.Range("J18").Value = 8.724621603287282
$ws.Range("M18").Value = 53.79585715174806

$ws.Range("B19").Value = 8.039700816234316
$ws.Range("D19").Value = 10.55226568403957
$ws.Range("E19").Value = 11.98538350331799
$ws.Range("F19").Value = 67.60049105773292
$ws.Range("G19").Value = 3.78968795510951
$ws.Range("J19").Value = 8.723648073698536
$ws.Range("M19").Value = 53.65017513428566

$ws.Range("B20").Value = 8.071938861875456
$ws.Range("D20").Value = 10.68177699966172
$ws.Range("E20").Value = 12.17518963629113
$ws.Range("F20").Value = 67.63482805033273
$ws.Range("G20").Value = 3.784853458939987
$ws.Range("J20").Value = 8.728162498652244
$ws.Range("M20").Value = 54.30332257957413

$ws.Range("B21").Value = 8.181419190891106
$ws.Range("D21").Value = 11.11521860516465
$ws.Range("E21").Value = 12.79382041798562
$ws.Range("F21").Value = 67.821516806052
$ws.Range("G21").Value = 3.768954292122657
$ws.Range("J21").Value = 8.745678937881603
$ws.Range("M21").Value = 56.45390311158878

$ws.Range("B22").Value = 8.253731562636954
$ws.Range("D22").Value = 11.39683388045961
$ws.Range("E22").Value = 13.18379973238041
$ws.Range("F22").Value = 67.998290440663
$ws.Range("G22").Value = 3.758812760373657
$ws.Range("J22").Value = 8.758894112083208
$ws.Range("M22").Value = 57.82533892729597

$ws.Range("B23").Value = 8.215074383878344
$ws.Range("D23").Value = 11.24670659360477
$ws.Range("E23").Value = 12.97697696319774
$ws.Range("F23").Value = 67.89883145606517
$ws.Range("G23").Value = 3.764203114842776
$ws.Range("J23").Value = 8.751678131490131
$ws.Range("M23").Value = 57.09657804308406

$ws.Range("B24").Value = 8.070163845706439
$ws.Range("D24").Value = 10.67467016436119
$ws.Range("E24").Value = 12.16483789435767
$ws.Range("F24").Value = 67.63267829137625
$ws.Range("G24").Value = 3.785117634630756
$ws.Range("J24").Value = 8.727905794292205
$ws.Range("M24").Value = 54.26761622332208

$ws.Range("B25").Value = 7.917816554385761
$ws.Range("D25").Value = 10.05320084181762
$ws.Range("E25").Value = 11.22806398573369
$ws.Range("F25").Value = 67.57013108688409
$ws.Range("G25").Value = 3.808792905596137
$ws.Range("J25").Value = 8.709743454426805
$ws.Range("M25").Value = 51.07921433467657
